{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, etc.)\n// in bold + dark slate color (#2C3E50) across specific resume bullet points,\n// matching the \"Implement quantitative metrics highlighting\" commit.\n//\n// Strategy: for each target paragraph (identified by a unique substring of\n// its original text), locate the metric token(s) inside that paragraph via\n// Paragraph.search() (scoped to the paragraph, so identical tokens that\n// appear in OTHER paragraphs - e.g. the professional summary or the key\n// projects section - are left untouched) and apply bold + color to just\n// that sub-range. Word automatically splits the run so the surrounding\n// text keeps its original (unbolded) formatting.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Each entry: a substring that uniquely identifies the target paragraph,\n// plus the ordered list of metric tokens inside it that should become\n// bold + colored (#2C3E50).\nconst edits = [\n  {\n    contains: \"\u2022 Discovered systematic race coding errors\",\n    tokens: [\"23%\", \"64%\"],\n  },\n  {\n    contains: \"\u2022 Utilized advanced sampling methods\",\n    tokens: [\"\u00b14.2%\", \"\u00b12.1%\", \"71%\", \"87%\"],\n  },\n  {\n    contains: \"\u2022 Trigonometric algorithm for boundary estimation\",\n    tokens: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    contains: \"\u2022 Built real-time FEC analysis systems\",\n    tokens: [\"$2\"],\n  },\n  {\n    contains: \"\u2022 Modernized legacy ETL processes\",\n    tokens: [\"57%\"],\n  },\n  {\n    contains: \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation\",\n    tokens: [\"73.5%\"],\n  },\n  {\n    contains: \"\u2022 $4.7M savings enabled nonprofit access\",\n    tokens: [\"$4.7M\"],\n  },\n  {\n    contains: \"\u2022 178% accuracy improvement in racial classification algorithms\",\n    tokens: [\"178%\"],\n  },\n];\n\nfor (const edit of edits) {\n  const target = paragraphs.items.find((p) => p.text.indexOf(edit.contains) !== -1);\n  if (!target) {\n    continue;\n  }\n  for (const token of edit.tokens) {\n    const results = target.search(token, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    for (let i = 0; i < results.items.length; i++) {\n      results.items[i].font.bold = true;\n      results.items[i].font.color = \"#2C3E50\";\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, etc.)\n# in bold + dark slate color (#2C3E50) across specific resume bullet points,\n# matching the \"Implement quantitative metrics highlighting\" commit.\n#\n# Strategy: for each target paragraph (identified by a unique substring of\n# its original text), locate the metric token(s) inside that paragraph's\n# Range via Range.Find.Execute (scoped to the paragraph, so identical\n# tokens appearing in OTHER paragraphs - e.g. the professional summary or\n# the key projects section - are left untouched) and apply Bold + Color to\n# just that found sub-range. Word automatically splits the run so the\n# surrounding text keeps its original (unbolded) formatting.\n\n$d = $word.ActiveDocument\n\n# #2C3E50 expressed as a Word/VBA BGR color integer (B*65536 + G*256 + R).\n$highlightColor = 0x50 * 65536 + 0x3E * 256 + 0x2C\n\n$edits = @(\n    @{ Contains = '\u2022 Discovered systematic race coding errors'; Tokens = @('23%', '64%') },\n    @{ Contains = '\u2022 Utilized advanced sampling methods'; Tokens = @('\u00b14.2%', '\u00b12.1%', '71%', '87%') },\n    @{ Contains = '\u2022 Trigonometric algorithm for boundary estimation'; Tokens = @('73.5%', '$4.7M') },\n    @{ Contains = '\u2022 Built real-time FEC analysis systems'; Tokens = @('$2') },\n    @{ Contains = '\u2022 Modernized legacy ETL processes'; Tokens = @('57%') },\n    @{ Contains = '\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation'; Tokens = @('73.5%') },\n    @{ Contains = '\u2022 $4.7M savings enabled nonprofit access'; Tokens = @('$4.7M') },\n    @{ Contains = '\u2022 178% accuracy improvement in racial classification algorithms'; Tokens = @('178%') }\n)\n\nforeach ($edit in $edits) {\n    foreach ($para in $d.Paragraphs) {\n        if ($para.Range.Text.Contains($edit.Contains)) {\n            foreach ($token in $edit.Tokens) {\n                $searchRange = $para.Range\n                $found = $searchRange.Find.Execute($token)\n                if ($found) {\n                    $searchRange.Font.Bold = 1\n                    $searchRange.Font.Color = $highlightColor\n                }\n            }\n            break\n        }\n    }\n}\n"}
